$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D29").Value = "개발팀 인턴의 좌충우돌 Figma i18n 플러그인 제작기"
$ws.Range("E29").Value = "https://blog.promedius.ai/figma-i18-plugin/"

$ws.Range("D43").Value = "jupyter notebook 브라우저 크롬으로 지정 2가지 방법"
$ws.Range("E43").Value = "https://nittaku.tistory.com/511"

$ws.Range("D50").Value = "Retraction of Nature paper puts Majorana research on a new path"
$ws.Range("E50").Value = "http://incredible.egloos.com/7515816"

$ws.Range("D51").Value = "[python] selenium으로 크롤링하는 중에 click() 안 먹힐 때"
$ws.Range("E51").Value = "https://bskyvision.com/1179"
